# Update countries & provincias Spain
#
# 1. Update the "last updated" timestamp banner.
# 2. Swap Grecia/Croacia: Croacia's numbers overtook Grecia's, so Croacia
#    (previously row 92) moves to row 91 and Grecia (previously row 91)
#    moves to row 92. Grecia's stats are unchanged; Croacia's stats are
#    refreshed.
# 3. Refresh the daily case counters for a handful of other countries.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Banner timestamp -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 2 de Septiembre de 2020 a las 10:30"

# --- Grecia / Croacia swap + refreshed data ---------------------------
$ws.Range("A91").Value = "Croacia"
$ws.Range("B91").Value = 10725
$ws.Range("C91").Value = 311
$ws.Range("D91").Value = 7968
$ws.Range("E91").Value = 2566
$ws.Range("F91").Value = 0
$ws.Range("G91").Value = 4
$ws.Range("H91").Value = 191

$ws.Range("A92").Value = "Grecia"
$ws.Range("B92").Value = 10524
$ws.Range("C92").Value = 0
$ws.Range("D92").Value = 3804
$ws.Range("E92").Value = 6449
$ws.Range("F92").Value = 0
$ws.Range("G92").Value = 0
$ws.Range("H92").Value = 271

# --- Other refreshed rows ----------------------------------------------
# Row 7 - Rusia
$ws.Range("B7").Value = 1005000
$ws.Range("C7").Value = 4952
$ws.Range("D7").Value = 821169
$ws.Range("E7").Value = 166417
$ws.Range("G7").Value = 115
$ws.Range("H7").Value = 17414

# Row 25 - Filipinas
$ws.Range("B25").Value = 226440
$ws.Range("C25").Value = 2218
$ws.Range("D25").Value = 158610
$ws.Range("E25").Value = 64207
$ws.Range("G25").Value = 27
$ws.Range("H25").Value = 3623

# Row 48
$ws.Range("D48").Value = 47865
$ws.Range("E48").Value = 17999

# Row 52
$ws.Range("B52").Value = 56901
$ws.Range("C52").Value = 49
$ws.Range("E52").Value = 1125

# Row 65
$ws.Range("B65").Value = 38243
$ws.Range("C65").Value = 47
$ws.Range("D65").Value = 29315
$ws.Range("E65").Value = 7519
$ws.Range("G65").Value = 3
$ws.Range("H65").Value = 1409

# Row 66
$ws.Range("D66").Value = 26189
$ws.Range("E66").Value = 10011

# Row 112
$ws.Range("E112").Value = 351
$ws.Range("G112").Value = 2
$ws.Range("H112").Value = 92

# Row 127
$ws.Range("D127").Value = 2883
$ws.Range("E127").Value = 197

# Row 130
$ws.Range("B130").Value = 2958
$ws.Range("C130").Value = 29
$ws.Range("D130").Value = 1874
$ws.Range("E130").Value = 998

# Row 137
$ws.Range("B137").Value = 2415
$ws.Range("C137").Value = 20
$ws.Range("D137").Value = 2130
$ws.Range("E137").Value = 221

# Row 150
$ws.Range("E150").Value = 824
$ws.Range("G150").Value = 1
$ws.Range("H150").Value = 10

# Row 155
$ws.Range("B155").Value = 1406
$ws.Range("C155").Value = 2
$ws.Range("E155").Value = 199

# Row 174
$ws.Range("B174").Value = 489
$ws.Range("C174").Value = 1
$ws.Range("D174").Value = 471
$ws.Range("E174").Value = 11
